# The upstream edit swaps the two 4-row record blocks A2:AY5 and A6:AY9 in
# place: every field of rows 2-5 ends up on rows 6-9 and vice versa (the
# record "Ids" in column A make this obvious: 75200047/75200049/75200051/
# 91846294 move from rows 2-5 down to rows 6-9, while 98677520/98677522/
# 99109084/99108983 move from rows 6-9 up to rows 2-5). Every other column
# just rides along with its row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$topRange = "A2:AY5"
$bottomRange = "A6:AY9"

$top = $ws.Range($topRange).Value()
$bottom = $ws.Range($bottomRange).Value()

# Plain ISO date strings (e.g. "2018-12-06") stored as text get silently
# re-interpreted as real dates by Excel when they're written back through
# .Value, which would change their type/format. Re-quote any such string
# with a leading apostrophe so it round-trips as text, exactly like it was
# before the swap.
$rows = $top.GetLength(0)
$cols = $top.GetLength(1)
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $v = $top[$r, $c]
        if ($v -is [string] -and $v -match '^\d{4}-\d{2}-\d{2}$') {
            $top[$r, $c] = "'" + $v
        }
        $v2 = $bottom[$r, $c]
        if ($v2 -is [string] -and $v2 -match '^\d{4}-\d{2}-\d{2}$') {
            $bottom[$r, $c] = "'" + $v2
        }
    }
}

$ws.Range($topRange).Value = $bottom
$ws.Range($bottomRange).Value = $top
